# Replace the 25 division-problem answers in the single table, cell by
# cell, using positional (row, col) addressing rather than text search —
# several of the old values repeat (e.g. "43÷8=5, 3" occurs twice) and
# map to different new values, so a global Find/Replace would be wrong.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "67÷4=16, 3"
$t.Cell(1, 2).Range.Text = "37÷2=18, 1"
$t.Cell(1, 3).Range.Text = "74÷9=8, 2"
$t.Cell(1, 4).Range.Text = "82÷8=10, 2"
$t.Cell(1, 5).Range.Text = "41÷6=6, 5"
$t.Cell(5, 1).Range.Text = "16÷3=5, 1"
$t.Cell(5, 2).Range.Text = "19÷7=2, 5"
$t.Cell(5, 3).Range.Text = "48÷8=6, 0"
$t.Cell(5, 4).Range.Text = "42÷3=14, 0"
$t.Cell(5, 5).Range.Text = "15÷7=2, 1"
$t.Cell(9, 1).Range.Text = "21÷8=2, 5"
$t.Cell(9, 2).Range.Text = "99÷8=12, 3"
$t.Cell(9, 3).Range.Text = "44÷6=7, 2"
$t.Cell(9, 4).Range.Text = "34÷8=4, 2"
$t.Cell(9, 5).Range.Text = "12÷7=1, 5"
$t.Cell(13, 1).Range.Text = "93÷5=18, 3"
$t.Cell(13, 2).Range.Text = "64÷2=32, 0"
$t.Cell(13, 3).Range.Text = "53÷6=8, 5"
$t.Cell(13, 4).Range.Text = "80÷5=16, 0"
$t.Cell(13, 5).Range.Text = "34÷9=3, 7"
$t.Cell(17, 1).Range.Text = "67÷7=9, 4"
$t.Cell(17, 2).Range.Text = "26÷6=4, 2"
$t.Cell(17, 3).Range.Text = "89÷3=29, 2"
$t.Cell(17, 4).Range.Text = "17÷8=2, 1"
$t.Cell(17, 5).Range.Text = "78÷8=9, 6"
